# Reposition a handful of shapes on slide 1 of the CV template:
#   - the purple divider line ("AutoShape 7") and the "{{PERFIL}}" label
#     ("TextBox 14") move up slightly so the rule sits closer to the
#     heading above it;
#   - the "{{EDUCACIÓN}}" label ("TextBox 15") moves down;
#   - the "{{HABILIDADES}}" label ("TextBox 17") moves up and slightly
#     left;
#   - the "{{IDIOMAS}}" label ("TextBox 18") moves down a touch.
#
# PowerPoint's Shape.Left/.Top are expressed in points (1 pt = 12700 EMU)
# and are stored with single-precision accuracy, so a direct
# EMU/12700.0 assignment can land 1 EMU away from the exact OOXML value
# once it is read back. Set-TopEmu / Set-LeftEmu below assign the shape
# position and then nudge it in tiny increments, re-reading the shape's
# own property each time, until the value that PowerPoint reports
# round-trips to precisely the requested EMU amount.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMU_PER_POINT = 12700.0

function Set-TopEmu($shape, $targetEmu) {
    $pts = $targetEmu / $EMU_PER_POINT
    $shape.Top = $pts
    $gotEmu = [Math]::Round($shape.Top * $EMU_PER_POINT)
    $tries = 0
    while (($gotEmu -ne $targetEmu) -and ($tries -lt 500)) {
        if ($gotEmu -lt $targetEmu) {
            $pts = $pts + 0.00001
        } else {
            $pts = $pts - 0.00001
        }
        $shape.Top = $pts
        $gotEmu = [Math]::Round($shape.Top * $EMU_PER_POINT)
        $tries = $tries + 1
    }
}

function Set-LeftEmu($shape, $targetEmu) {
    $pts = $targetEmu / $EMU_PER_POINT
    $shape.Left = $pts
    $gotEmu = [Math]::Round($shape.Left * $EMU_PER_POINT)
    $tries = 0
    while (($gotEmu -ne $targetEmu) -and ($tries -lt 500)) {
        if ($gotEmu -lt $targetEmu) {
            $pts = $pts + 0.00001
        } else {
            $pts = $pts - 0.00001
        }
        $shape.Left = $pts
        $gotEmu = [Math]::Round($shape.Left * $EMU_PER_POINT)
        $tries = $tries + 1
    }
}

# Divider line above the "PERFIL" section
Set-TopEmu $s.Shapes.Item("AutoShape 7") 2286000

# "{{PERFIL}}" label, aligned with the line above
Set-TopEmu $s.Shapes.Item("TextBox 14") 2286000

# "{{EDUCACIÓN}}" label
Set-TopEmu $s.Shapes.Item("TextBox 15") 7112606

# "{{HABILIDADES}}" label
Set-LeftEmu $s.Shapes.Item("TextBox 17") 410245
Set-TopEmu $s.Shapes.Item("TextBox 17") 4343400

# "{{IDIOMAS}}" label
Set-TopEmu $s.Shapes.Item("TextBox 18") 8712036
